$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "kvw5270"
$ws.Range("B2").Value = "04/14/2020 00:51:46"
$ws.Range("C2").Value = "04/14/2020 00:51:59"
